$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.1485482781904119
$ws.Cells.Item(2, 2).Value = 1689843629.714195
$ws.Cells.Item(2, 3).Value = 0.3743307501007689
$ws.Cells.Item(2, 4).Value = 1689843629.799342
$ws.Cells.Item(2, 5).Value = 0.08514714241027832

$ws.Cells.Item(3, 1).Value = 0.08912896691424713
$ws.Cells.Item(3, 2).Value = 1689843634.716197
$ws.Cells.Item(3, 3).Value = 0.06667587402748519
$ws.Cells.Item(3, 4).Value = 1689843634.800273
$ws.Cells.Item(3, 5).Value = 0.0840761661529541

$ws.Cells.Item(4, 1).Value = 0.03308575286968264
$ws.Cells.Item(4, 2).Value = 1689843639.701746
$ws.Cells.Item(4, 3).Value = 0.7917310801664017
$ws.Cells.Item(4, 4).Value = 1689843639.808546
$ws.Cells.Item(4, 5).Value = 0.1067993640899658

$ws.Cells.Item(5, 1).Value = 0.224510465901418
$ws.Cells.Item(5, 2).Value = 1689843644.717801
$ws.Cells.Item(5, 3).Value = 0.5067153602202465
$ws.Cells.Item(5, 4).Value = 1689843644.799302
$ws.Cells.Item(5, 5).Value = 0.08150076866149902

$ws.Cells.Item(6, 1).Value = 0.1583389601620526
$ws.Cells.Item(6, 2).Value = 1689843649.719811
$ws.Cells.Item(6, 3).Value = 0.2526672193111211
$ws.Cells.Item(6, 4).Value = 1689843649.798629
$ws.Cells.Item(6, 5).Value = 0.07881784439086914

$ws.Cells.Item(7, 1).Value = 0.1056718433490885
$ws.Cells.Item(7, 2).Value = 1689843654.715734
$ws.Cells.Item(7, 3).Value = 0.01247903626919362
$ws.Cells.Item(7, 4).Value = 1689843653.838972
$ws.Cells.Item(7, 5).Value = 0.8767616748809814

$ws.Cells.Item(8, 1).Value = 0.0462525320729237
$ws.Cells.Item(8, 2).Value = 1689843659.718307
$ws.Cells.Item(8, 3).Value = 0.01249728874658717
$ws.Cells.Item(8, 4).Value = 1689843659.248995
$ws.Cells.Item(8, 5).Value = 0.4693119525909424

$ws.Cells.Item(9, 1).Value = 0.2045914922349764
$ws.Cells.Item(9, 2).Value = 1689843664.715051
$ws.Cells.Item(9, 3).Value = 0.4027382366585798
$ws.Cells.Item(9, 4).Value = 1689843664.798724
$ws.Cells.Item(9, 5).Value = 0.08367323875427246

$ws.Cells.Item(10, 1).Value = 0.1451721809588116
$ws.Cells.Item(10, 2).Value = 1689843669.713756
$ws.Cells.Item(10, 3).Value = 0.01267664975777441
$ws.Cells.Item(10, 4).Value = 1689843668.979332
$ws.Cells.Item(10, 5).Value = 0.7344245910644531

$ws.Cells.Item(11, 1).Value = 0.08575286968264685
$ws.Cells.Item(11, 2).Value = 1689843674.696821
$ws.Cells.Item(11, 3).Value = 0.0125041638464054
$ws.Cells.Item(11, 4).Value = 1689843672.778796
$ws.Cells.Item(11, 5).Value = 1.918024778366089

$ws.Cells.Item(12, 1).Value = 0.08575286968264685
$ws.Cells.Item(12, 2).Value = 1689843674.696821
$ws.Cells.Item(12, 3).Value = 0.02056579637840428
$ws.Cells.Item(12, 4).Value = 1689843675.785807
$ws.Cells.Item(12, 5).Value = 1.088986396789551

$ws.Cells.Item(13, 1).Value = 0.02295746117488183
$ws.Cells.Item(13, 2).Value = 1689843679.6954
$ws.Cells.Item(13, 3).Value = 0.5741404984447368
$ws.Cells.Item(13, 4).Value = 1689843679.800069
$ws.Cells.Item(13, 5).Value = 0.1046686172485352

